# Refresh the cryptocurrency price/volume snapshot in Sheet1.
#
# Column D ("Price") holds plain-text numbers (some of them look like
# thousand-grouped numbers, e.g. "45.415.69") and column E ("Volume(1h)")
# holds padded percentage strings, e.g. "  -0.97%  ". Both must stay TEXT,
# exactly as authored upstream - not get auto-coerced into Excel numbers.
#
# Rows 36-40 also have their Coin/Link (columns B/C) values reshuffled as
# part of this refresh, in addition to their Price/Volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '45.415.69'
$ws.Range('E2').Value = '  -0.97%  '
# Row 3
$ws.Range('D3').Value = '2.374.51'
$ws.Range('E3').Value = '  -2.09%  '
# Row 5
$ws.Range('D5').Value = '''320.36'
$ws.Range('E5').Value = '  +0.11%  '
# Row 6
$ws.Range('D6').Value = '''108.84'
$ws.Range('E6').Value = '  -5.97%  '
# Row 7
$ws.Range('D7').Value = '''0.638'
$ws.Range('E7').Value = '  +0.50%  '
# Row 8
$ws.Range('E8').Value = '  +0.07%  '
# Row 9
$ws.Range('E9').Value = '  -2.63%  '
# Row 10
$ws.Range('D10').Value = '''41.01'
$ws.Range('E10').Value = '  -4.98%  '
# Row 11
$ws.Range('E11').Value = '  -2.28%  '
# Row 12
$ws.Range('D12').Value = '''8.52'
$ws.Range('E12').Value = '  -2.28%  '
# Row 13
$ws.Range('E13').Value = '  +0.21%  '
# Row 14
$ws.Range('D14').Value = '''0.982'
$ws.Range('E14').Value = '  -4.68%  '
# Row 15
$ws.Range('D15').Value = '2.735.50'
$ws.Range('E15').Value = '  -1.99%  '
# Row 16
$ws.Range('D16').Value = '''15.46'
$ws.Range('E16').Value = '  -3.59%  '
# Row 17
$ws.Range('D17').Value = '2.382.12'
$ws.Range('E17').Value = '  -1.48%  '
# Row 18
$ws.Range('D18').Value = '45.363.30'
$ws.Range('E18').Value = '  -1.07%  '
# Row 19
$ws.Range('D19').Value = '''15.14'
$ws.Range('E19').Value = '  +12.79%  '
# Row 20
$ws.Range('D20').Value = '''7.34'
$ws.Range('E20').Value = '  -4.05%  '
# Row 21
$ws.Range('D21').Value = '''0.0000107'
$ws.Range('E21').Value = '  -2.52%  '
# Row 22
$ws.Range('D22').Value = '''3.67'
$ws.Range('E22').Value = '  +3.21%  '
# Row 23
$ws.Range('D23').Value = '''73.33'
$ws.Range('E23').Value = '  -2.60%  '
# Row 24
$ws.Range('D24').Value = '''264.41'
$ws.Range('E24').Value = '  -2.06%  '
# Row 25
$ws.Range('E25').Value = '  -1.49%  '
# Row 26
$ws.Range('E26').Value = '  +0.18%  '
# Row 27
$ws.Range('E27').Value = '  -1.98%  '
# Row 28
$ws.Range('D28').Value = '''11.26'
$ws.Range('E28').Value = '  -1.20%  '
# Row 29
$ws.Range('E29').Value = '  -1.83%  '
# Row 30
$ws.Range('E30').Value = '  -2.90%  '
# Row 31
$ws.Range('D31').Value = '''0.0952'
$ws.Range('E31').Value = '  -1.33%  '
# Row 32
$ws.Range('D32').Value = '''37.28'
$ws.Range('E32').Value = '  -6.39%  '
# Row 33
$ws.Range('D33').Value = '''168.62'
$ws.Range('E33').Value = '  -2.98%  '
# Row 34
$ws.Range('E34').Value = '  -4.14%  '
# Row 35
$ws.Range('E35').Value = '  -0.41%  '
# Row 36
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '''3.28'
$ws.Range('E36').Value = '  +3.82%  '
# Row 37
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '''0.118'
$ws.Range('E37').Value = '  -3.93%  '
# Row 38
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = '''4.73'
$ws.Range('E38').Value = '  -5.37%  '
# Row 39
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''1.95'
$ws.Range('E39').Value = '  +6.47%  '
# Row 40
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').Value = '''4.03'
$ws.Range('E40').Value = '  -4.05%  '
# Row 41
$ws.Range('E41').Value = '  -3.18%  '
# Row 42
$ws.Range('D42').Value = '''98.41'
$ws.Range('E42').Value = '  -4.22%  '
# Row 43
$ws.Range('D43').Value = '''70.32'
$ws.Range('E43').Value = '  -2.75%  '
# Row 44
$ws.Range('D44').Value = '1.879.06'
$ws.Range('E44').Value = '  +12.82%  '
# Row 45
$ws.Range('E45').Value = '  -4.42%  '
# Row 46
$ws.Range('E46').Value = '  -5.15%  '
# Row 47
$ws.Range('D47').Value = '''6.04'
$ws.Range('E47').Value = '  +2.61%  '
# Row 48
$ws.Range('E48').Value = '  +0.14%  '
# Row 49
$ws.Range('D49').Value = '''84.63'
$ws.Range('E49').Value = '  +6.20%  '
# Row 50
$ws.Range('D50').Value = '''112.60'
$ws.Range('E50').Value = '  -4.44%  '
# Row 51
$ws.Range('D51').Value = '''9.36'
$ws.Range('E51').Value = '  -1.26%  '
